$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds data rows 2..129 (row 1 is the header).
# A brand-new record is inserted as the new row 39, pushing the former
# rows 39..129 down by one (to 40..130).

$ws.Rows.Item(39).Insert()

$ws.Cells.Item(39, 1).Value  = 5
$ws.Cells.Item(39, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(39, 3).Value  = "Maule"
$ws.Cells.Item(39, 4).Value  = 44533
$ws.Cells.Item(39, 5).Value  = 7
$ws.Cells.Item(39, 6).Value  = 100112024
$ws.Cells.Item(39, 7).Value  = "Choclo"
$ws.Cells.Item(39, 8).Value  = "Choclero"
$ws.Cells.Item(39, 9).Value  = "Primera"
$ws.Cells.Item(39, 10).Value = 10000
$ws.Cells.Item(39, 11).Value = 400
$ws.Cells.Item(39, 12).Value = 400
$ws.Cells.Item(39, 13).Value = 400
$ws.Cells.Item(39, 14).Value = "$/unidad"
$ws.Cells.Item(39, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(39, 16).Value = 400
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"
